$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to write a value as plain text (avoids Excel auto-converting
# numeric-looking strings like "628.53" into floating point numbers),
# while leaving the cell's style/format untouched in the saved file.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "68.905.47"
Set-TextValue "E2" "  +1.48%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.773.22"
Set-TextValue "E3" "  -0.31%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  -0.08%  "

# Row 5 - BNB
Set-TextValue "D5" "628.53"
Set-TextValue "E5" "  +4.10%  "

# Row 6 - Solana
Set-TextValue "D6" "165.00"
Set-TextValue "E6" "  +0.87%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.771.01"
Set-TextValue "E7" "  -0.28%  "

# Row 8 - USDC
Set-TextValue "E8" "  -0.11%  "

# Row 9 - XRP
Set-TextValue "E9" "  +1.35%  "

# Row 10 - Dogecoin
Set-TextValue "E10" "  +0.87%  "

# Row 11 - Cardano
Set-TextValue "E11" "  +2.34%  "

# Row 12 - Toncoin
Set-TextValue "E12" "  +0.01%  "

# Row 13 - ShibaInu
Set-TextValue "E13" "  -0.97%  "

# Row 14 - Avalanche
Set-TextValue "D14" "35.14"
Set-TextValue "E14" "  +0.17%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "4.410.10"
Set-TextValue "E15" "  -0.14%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.800.05"
Set-TextValue "E16" "  +0.72%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "68.962.66"
Set-TextValue "E17" "  +1.61%  "

# Row 18 - Chainlink
Set-TextValue "D18" "17.57"
Set-TextValue "E18" "  -3.36%  "

# Row 20 - Polkadot
Set-TextValue "D20" "7.02"
Set-TextValue "E20" "  +0.34%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "465.84"
Set-TextValue "E21" "  +1.56%  "

# Row 22 - Uniswap
Set-TextValue "E22" "  +0.71%  "

# Row 23 - Polygon
Set-TextValue "D23" "0.705"
Set-TextValue "E23" "  +2.08%  "

# Row 24 - Litecoin
Set-TextValue "D24" "82.95"
Set-TextValue "E24" "  -0.24%  "

# Row 25 - PEPE
Set-TextValue "E25" "  +0.03%  "

# Row 26 - InternetComputer(DFINITY)
Set-TextValue "D26" "11.98"
Set-TextValue "E26" "  +0.86%  "

# Row 27 - Fetch.AI
Set-TextValue "E27" "  +2.96%  "

# Row 28 - RenderToken
Set-TextValue "D28" "10.01"
Set-TextValue "E28" "  +1.17%  "

# Row 30 - WrappedeETH (only price changes)
Set-TextValue "D30" "3.924.52"

# Row 31 - PancakeSwap
Set-TextValue "E31" "  +2.75%  "

# Row 32 - ImmutableX
Set-TextValue "E32" "  +2.11%  "

# Row 33 - NEARProtocol
Set-TextValue "D33" "7.11"
Set-TextValue "E33" "  -1.37%  "

# Row 34 - EthereumClassic
Set-TextValue "D34" "28.69"
Set-TextValue "E34" "  -1.08%  "

# Row 35 - Kaspa
Set-TextValue "D35" "0.171"
Set-TextValue "E35" "  +16.04%  "

# Row 37 - RenzoRestakedETH
Set-TextValue "D37" "3.725.87"
Set-TextValue "E37" "  -0.24%  "

# Row 38 - Aptos
Set-TextValue "D38" "8.92"
Set-TextValue "E38" "  +0.04%  "

# Row 39 - Hedera
Set-TextValue "E39" "  +1.80%  "

# Row 40 - dogwifhat
Set-TextValue "E40" "  +2.65%  "

# Row 41 - Filecoin
Set-TextValue "D41" "5.80"
Set-TextValue "E41" "  +0.01%  "

# Row 42 - Mantle
Set-TextValue "D42" "0.968"
Set-TextValue "E42" "  -1.11%  "

# Row 43 - FirstDigitalUSD
Set-TextValue "E43" "  +0.06%  "

# Row 45 - Monero
Set-TextValue "D45" "155.50"
Set-TextValue "E45" "  +1.96%  "

# Row 46 - Arweave
Set-TextValue "D46" "42.98"
Set-TextValue "E46" "  -1.27%  "

# Row 47 - TheGraph
Set-TextValue "D47" "0.294"
Set-TextValue "E47" "  +0.24%  "

# Row 48 - OKB
Set-TextValue "D48" "46.62"
Set-TextValue "E48" "  -0.91%  "

# Row 49 - Stacks
Set-TextValue "D49" "1.89"
Set-TextValue "E49" "  +3.59%  "

# Row 50 - Cosmos
Set-TextValue "D50" "8.35"
Set-TextValue "E50" "  +0.61%  "

# Row 51 - ONDO
Set-TextValue "E51" "  -1.07%  "
